# JS08 - Refactor: fitur import data supplier
#
# The "Supplier" sheet drops its leading "supplier_id" column. The
# remaining columns (supplier_kode, supplier_nama, supplier_alamat)
# shift one column to the left (B->A, C->B, D->C) while keeping all
# their original values/pairings intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Used range is A1:D4 - 4 header/data rows, 4 columns.
$lastRow = 4

for ($r = 1; $r -le $lastRow; $r++) {
    $kode   = $ws.Cells.Item($r, 2).Value2
    $nama   = $ws.Cells.Item($r, 3).Value2
    $alamat = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 1).Value2 = $kode
    $ws.Cells.Item($r, 2).Value2 = $nama
    $ws.Cells.Item($r, 3).Value2 = $alamat
}

# Drop the now-stale last column (old "supplier_alamat" data / header)
# entirely so the sheet's used range shrinks back down to A1:C4.
$ws.Range("D1:D4").Clear()

# Leave the selection where the edit ended up.
$ws.Range("C9").Select()
